$wb = $excel.ActiveWorkbook

# --- Sheet 1: Overview ---
# Status text for zh-cn / de-de rows moves from "Ready for handoff" to "In Translation"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

# Columns E and F auto-narrow to fit the shorter "In Translation" text
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

# --- Sheet 2: zh-cn ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Columns.Item(3).ColumnWidth = 12.5

# --- Sheet 3: de-de ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Columns.Item(3).ColumnWidth = 12.5
